$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates (Volume/Number and week-covering dates) ---
# Keep run-level font formatting implicit via the cell style; only the
# textual content changes (all runs in these cells already share the same
# font/size/color, so a plain text replace is visually identical).
$c8 = $ws.Range("A8")
$n8 = $c8.Characters().Count
$c8.Characters(1, $n8).Text = "Volume 29   Number  45"

$c9 = $ws.Range("C9")
$n9 = $c9.Characters().Count
$c9.Characters(1, $n9).Text = "Report Covering the Week  11/7/2022  Through  11/13/2022"

# --- Cells that change type (text placeholder "0"/"***.*"  <->  numeric) ---
# Use PasteSpecial(Formats) from a donor cell of the desired style so the
# cell lands on the same shared numeric/text style used elsewhere in the
# table, then set the value.

# C15: "0" (text) -> 1 (number)
$ws.Range("F15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = 1

# C16: "0" (text) -> 2 (number)
$ws.Range("F15").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Value = 2

# D16: 3 (number) -> "0" (text)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = "'0"

# E16: -100 (number) -> "***.*" (text)
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = "'***.*"

# C26: "0" (text) -> 1 (number)
$ws.Range("F15").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null
$ws.Range("C26").Value = 1

# F28: 1 (number) -> "0" (text)
$ws.Range("G28").Copy() | Out-Null
$ws.Range("F28").PasteSpecial(-4122) | Out-Null
$ws.Range("F28").Value = "'0"

# F29: 1 (number) -> "0" (text)
$ws.Range("G29").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122) | Out-Null
$ws.Range("F29").Value = "'0"

# --- Plain numeric updates across the weekly crime-stat table ---
$ws.Range("N14").Value = -90
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 19
$ws.Range("K15").Value = 111.111111111111
$ws.Range("L15").Value = 26.666666666666
$ws.Range("M15").Value = 137.5
$ws.Range("N15").Value = 5.555555555555
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 42.857142857142
$ws.Range("I16").Value = 81
$ws.Range("K16").Value = 22.727272727272
$ws.Range("L16").Value = -13.829787234042
$ws.Range("N16").Value = -87.223974763406
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 12
$ws.Range("I17").Value = 189
$ws.Range("J17").Value = 145
$ws.Range("K17").Value = 30.344827586206
$ws.Range("L17").Value = 45.384615384615
$ws.Range("M17").Value = 30.344827586206
$ws.Range("N17").Value = -32.978723404255
$ws.Range("C18").Value = 10
$ws.Range("E18").Value = 233.333333333333
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 18.75
$ws.Range("I18").Value = 167
$ws.Range("J18").Value = 142
$ws.Range("K18").Value = 17.605633802816
$ws.Range("L18").Value = 14.383561643835
$ws.Range("M18").Value = -49.698795180722
$ws.Range("N18").Value = -88.806970509383
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -18.181818181818
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 26.829268292682
$ws.Range("I19").Value = 546
$ws.Range("J19").Value = 388
$ws.Range("K19").Value = 40.721649484536
$ws.Range("L19").Value = 62.985074626865
$ws.Range("M19").Value = 66.463414634146
$ws.Range("N19").Value = -8.080808080808
$ws.Range("C20").Value = 4
$ws.Range("I20").Value = 108
$ws.Range("K20").Value = 74.193548387096
$ws.Range("L20").Value = -2.702702702702
$ws.Range("M20").Value = -11.475409836065
$ws.Range("N20").Value = -93.233082706766
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 68.75
$ws.Range("F21").Value = 105
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = 34.615384615384
$ws.Range("I21").Value = 1111
$ws.Range("J21").Value = 813
$ws.Range("K21").Value = 36.654366543665
$ws.Range("L21").Value = 33.213429256594
$ws.Range("M21").Value = 1.276207839562
$ws.Range("N21").Value = -75.983571119757
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -38.461538461538
$ws.Range("F24").Value = 70
$ws.Range("G24").Value = 84
$ws.Range("H24").Value = -16.666666666666
$ws.Range("I24").Value = 971
$ws.Range("J24").Value = 748
$ws.Range("K24").Value = 29.812834224598
$ws.Range("L24").Value = 28.439153439153
$ws.Range("M24").Value = 20.173267326732
$ws.Range("C25").Value = 6
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = -10
$ws.Range("I25").Value = 304
$ws.Range("J25").Value = 268
$ws.Range("K25").Value = 13.432835820895
$ws.Range("L25").Value = 7.801418439716
$ws.Range("M25").Value = -15.789473684210
$ws.Range("F26").Value = 6
$ws.Range("H26").Value = 200
$ws.Range("I26").Value = 26
$ws.Range("K26").Value = 85.714285714285
$ws.Range("L26").Value = 30
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 16.666666666666
$ws.Range("I27").Value = 71
$ws.Range("J27").Value = 63
$ws.Range("K27").Value = 12.698412698412
$ws.Range("L27").Value = 77.5
